$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 134-137: odds data refreshed (rows effectively reshuffled) ---
# Row 134 (id 132) takes on the data previously held by row 137
$ws.Range("B134").Value = 7483306
$ws.Range("F134").Value = "Tecnico Universitario"
$ws.Range("G134").Value = "Club Atletico Libertad"
$ws.Range("H134").Value = 1
$ws.Range("I134").Value = 1
$ws.Range("J134").Value = "D"
$ws.Range("K134").Value = 1.5
$ws.Range("L134").Value = 4.333
$ws.Range("M134").Value = 5.75
$ws.Range("N134").Value = 1.533
$ws.Range("O134").Value = 4.2
$ws.Range("P134").Value = 5.5
$ws.Range("Q134").Value = -1
$ws.Range("R134").Value = 1.925
$ws.Range("S134").Value = 1.875
$ws.Range("T134").Value = 2.25
$ws.Range("U134").Value = 1.8
$ws.Range("V134").Value = 2
$ws.Range("W134").Value = -1
$ws.Range("X134").Value = 3.2
$ws.Range("Y134").Value = -1
$ws.Range("Z134").Value = -1
$ws.Range("AA134").Value = 0.875
$ws.Range("AB134").Value = -0.5
$ws.Range("AC134").Value = 0.5

# Row 135 (id 133) takes on the data previously held by row 136
$ws.Range("B135").Value = 7482867
$ws.Range("F135").Value = "Cumbaya FC"
$ws.Range("G135").Value = "LDU Quito"
$ws.Range("H135").Value = 1
$ws.Range("I135").Value = 2
$ws.Range("J135").Value = "A"
$ws.Range("K135").Value = 5.25
$ws.Range("L135").Value = 3.75
$ws.Range("M135").Value = 1.65
$ws.Range("N135").Value = 9
$ws.Range("O135").Value = 4.5
$ws.Range("P135").Value = 1.363
$ws.Range("Q135").Value = 1.25
$ws.Range("R135").Value = 1.975
$ws.Range("S135").Value = 1.825
$ws.Range("T135").Value = 2.5
$ws.Range("U135").Value = 1.825
$ws.Range("V135").Value = 1.975
$ws.Range("W135").Value = -1
$ws.Range("X135").Value = -1
$ws.Range("Y135").Value = 0.363
$ws.Range("Z135").Value = 0.4875
$ws.Range("AA135").Value = -0.5
$ws.Range("AB135").Value = 0.825
$ws.Range("AC135").Value = -1

# Row 136 (id 134) takes on the data previously held by row 135
$ws.Range("B136").Value = 7483188
$ws.Range("F136").Value = "Gualaceo SC"
$ws.Range("G136").Value = "Emelec"
$ws.Range("H136").Value = 0
$ws.Range("I136").Value = 2
$ws.Range("J136").Value = "A"
$ws.Range("K136").Value = 3.6
$ws.Range("L136").Value = 3.3
$ws.Range("M136").Value = 2.05
$ws.Range("N136").Value = 2.6
$ws.Range("O136").Value = 3.25
$ws.Range("P136").Value = 2.75
$ws.Range("Q136").Value = 0
$ws.Range("R136").Value = 1.8
$ws.Range("S136").Value = 2
$ws.Range("T136").Value = 2.5
$ws.Range("U136").Value = 1.975
$ws.Range("V136").Value = 1.825
$ws.Range("W136").Value = -1
$ws.Range("X136").Value = -1
$ws.Range("Y136").Value = 1.75
$ws.Range("Z136").Value = -1
$ws.Range("AA136").Value = 1
$ws.Range("AB136").Value = -1
$ws.Range("AC136").Value = 0.825

# Row 137 (id 135) takes on the data previously held by row 134
$ws.Range("B137").Value = 7482832
$ws.Range("F137").Value = "Barcelona Guayaquil"
$ws.Range("G137").Value = "Guayaquil City"
$ws.Range("H137").Value = 2
$ws.Range("I137").Value = 1
$ws.Range("J137").Value = "H"
$ws.Range("K137").Value = 1.363
$ws.Range("L137").Value = 5
$ws.Range("M137").Value = 7.5
$ws.Range("N137").Value = 1.444
$ws.Range("O137").Value = 4
$ws.Range("P137").Value = 8
$ws.Range("Q137").Value = -1.25
$ws.Range("R137").Value = 2.05
$ws.Range("S137").Value = 1.75
$ws.Range("T137").Value = 2.5
$ws.Range("U137").Value = 1.95
$ws.Range("V137").Value = 1.85
$ws.Range("W137").Value = 0.444
$ws.Range("X137").Value = -1
$ws.Range("Y137").Value = -1
$ws.Range("Z137").Value = -0.5
$ws.Range("AA137").Value = 0.375
$ws.Range("AB137").Value = 0.95
$ws.Range("AC137").Value = -1

# --- Rows 192-197: odds values updated ---
$ws.Range("N192").Value = 1.444
$ws.Range("P192").Value = 6.5
$ws.Range("R192").Value = 1.95
$ws.Range("S192").Value = 1.85
$ws.Range("U192").Value = 1.9
$ws.Range("V192").Value = 1.9

$ws.Range("N193").Value = 3.2
$ws.Range("O193").Value = 3.6
$ws.Range("P193").Value = 2.05
$ws.Range("Q193").Value = 0.25
$ws.Range("R193").Value = 1.95
$ws.Range("S193").Value = 1.85
$ws.Range("U193").Value = 1.85
$ws.Range("V193").Value = 1.95

$ws.Range("P194").Value = 7
$ws.Range("R194").Value = 1.8
$ws.Range("S194").Value = 2

$ws.Range("R195").Value = 2.025
$ws.Range("S195").Value = 1.775
$ws.Range("U195").Value = 1.95
$ws.Range("V195").Value = 1.85

$ws.Range("U197").Value = 2
$ws.Range("V197").Value = 1.8
